$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.502.56'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.814.55'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''228.49'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('E6').Value = '  +4.22%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '''34.89'
$ws.Range('E8').Value = '  +6.63%  '
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').Value = '''0.0694'
$ws.Range('E10').Value = '  +0.41%  '
$ws.Range('D11').Value = '''0.0952'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').Value = '2.078.57'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '''11.32'
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('D14').Value = '1.817.52'
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').Value = '''0.646'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').Value = '34.520.24'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').Value = '''4.36'
$ws.Range('E17').Value = '  +2.56%  '
$ws.Range('D18').Value = '''69.35'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').Value = '0.0₃0799'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '''245.65'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').Value = '''11.54'
$ws.Range('E21').Value = '  +2.66%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').Value = '''171.87'
$ws.Range('E24').Value = '  +4.24%  '
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('D26').Value = '''7.54'
$ws.Range('E26').Value = '  +4.27%  '
$ws.Range('D27').Value = '''16.81'
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('E28').Value = '  +2.46%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '''4.01'
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('D33').Value = '''3.83'
$ws.Range('E33').Value = '  +1.14%  '
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('B35').Value = 'RenderToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').Value = '''2.57'
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.401.38'
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').Value = '''0.680'
$ws.Range('E37').Value = '  +1.50%  '
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '''2.86'
$ws.Range('E40').Value = '  +4.28%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '''83.12'
$ws.Range('E41').Value = '  -1.98%  '
$ws.Range('D42').Value = '''0.953'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').Value = '''13.86'
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('D46').Value = '''0.0507'
$ws.Range('E46').Value = '  -3.14%  '
$ws.Range('D47').Value = '''6.04'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').Value = '1.978.75'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').Value = '''105.71'
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('E51').Value = '  +0.19%  '
